$wb = $excel.ActiveWorkbook
$replacements = @(
    @('ADR', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)'),
    @('CSV', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)'),
    @('DP ', 'DP - Democratic Party (Demokratesch Partei, DP)'),
    @('G -', 'G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)'),
    @('LSA', 'LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('L -', 'L - The Left (Déi Lénk, L)'),
    @('Pir', 'Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)'),
    @('ADR+DP ', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)'),
    @('ADR+G -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)'),
    @('ADR+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('CSV+DP ', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)'),
    @('CSV+G -', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)'),
    @('CSV+LSA', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('DP +G -', 'DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)'),
    @('DP +LSA', 'DP - Democratic Party (Demokratesch Partei, DP)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('G -+LSA', 'G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+CSV+DP ', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)'),
    @('ADR+CSV+G -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)'),
    @('ADR+CSV+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+DP +G -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)'),
    @('ADR+DP +LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+G -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('CSV+DP +G -', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)'),
    @('CSV+DP +LSA', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('CSV+G -+LSA', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('DP +G -+LSA', 'DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+CSV+DP +G -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)'),
    @('ADR+CSV+DP +LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+CSV+G -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+DP +G -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('CSV+DP +G -+LSA', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+CSV+DP +G -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+L -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+L - The Left (Déi Lénk, L)'),
    @('CSV+L -', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+L - The Left (Déi Lénk, L)'),
    @('DP +L -', 'DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)'),
    @('G -+L -', 'G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)'),
    @('L -+LSA', 'L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+CSV+L -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+L - The Left (Déi Lénk, L)'),
    @('ADR+DP +L -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)'),
    @('ADR+G -+L -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)'),
    @('ADR+L -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('CSV+DP +L -', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)'),
    @('CSV+G -+L -', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)'),
    @('CSV+L -+LSA', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('DP +G -+L -', 'DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)'),
    @('DP +L -+LSA', 'DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('G -+L -+LSA', 'G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+CSV+DP +L -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)'),
    @('ADR+CSV+G -+L -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)'),
    @('ADR+CSV+L -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+DP +G -+L -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)'),
    @('ADR+DP +L -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+G -+L -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('CSV+DP +G -+L -', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)'),
    @('CSV+DP +L -+LSA', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('CSV+G -+L -+LSA', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('DP +G -+L -+LSA', 'DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+CSV+DP +G -+L -', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)'),
    @('ADR+CSV+DP +L -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+CSV+G -+L -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+DP +G -+L -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('CSV+DP +G -+L -+LSA', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+CSV+DP +G -+L -+LSA', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)'),
    @('ADR+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('DP +Pir', 'DP - Democratic Party (Demokratesch Partei, DP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('G -+Pir', 'G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('L -+Pir', 'L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('LSA+Pir', 'LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+DP +Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+G -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+L -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+DP +Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+G -+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+L -+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+LSA+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('DP +G -+Pir', 'DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('DP +L -+Pir', 'DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('DP +LSA+Pir', 'DP - Democratic Party (Demokratesch Partei, DP)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('G -+L -+Pir', 'G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('G -+LSA+Pir', 'G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('L -+LSA+Pir', 'L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+DP +Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+G -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+L -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+DP +G -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+DP +L -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+DP +LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+G -+L -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+G -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+L -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+DP +G -+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+DP +L -+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+DP +LSA+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+G -+L -+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+G -+LSA+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+L -+LSA+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('DP +G -+L -+Pir', 'DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('DP +G -+LSA+Pir', 'DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('DP +L -+LSA+Pir', 'DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('G -+L -+LSA+Pir', 'G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+DP +G -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+DP +L -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+DP +LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+G -+L -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+G -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+L -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+DP +G -+L -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+DP +G -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+DP +L -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+G -+L -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+DP +G -+L -+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+DP +G -+LSA+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+DP +L -+LSA+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+G -+L -+LSA+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('DP +G -+L -+LSA+Pir', 'DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+DP +G -+L -+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+DP +G -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+DP +L -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+G -+L -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+DP +G -+L -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('CSV+DP +G -+L -+LSA+Pir', 'CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('ADR+CSV+DP +G -+L -+LSA+Pir', 'ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)'),
    @('(''ADR+CSV'', ''DP +G -+L -+LSA'')', '(''ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)'', "DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)")'),
    @('(''DP +G -+L -+LSA'', ''ADR+CSV'')', '("DP - Democratic Party (Demokratesch Partei, DP)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)", ''ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)'')'),
    @('(''CSV+G -'', ''ADR+DP +L -+LSA+Pir'')', '(''CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)'', "ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)")'),
    @('(''ADR+DP +L -+LSA+Pir'', ''CSV+G -'')', '("ADR - Alternative Democratic Reform Party (Alternativ Demokratesch Reformpartei, ADR), known until 02 November 1992 as Action Committee 5/6ths Pensions for Everyone (ADR, Aktiounskomitee 5/6 Pensioun fir jiddfereen), known from 02 November 1992 until 02 April 2006 as Action Committee for Democracy and Pensions Justice (ADR, Aktiounskomitee fir Demokratie a  Rentegerechtegkeet)+DP - Democratic Party (Demokratesch Partei, DP)+L - The Left (Déi Lénk, L)+LSAP - Socialist Workers'' Party (Letzebuerger Sozialistesch Arbechterparteil, LSAP)+Pirat - Prirate Party (Piratepartei Lëtzebuerg, Pirat)", ''CSV - Christian Social Party (Chrestlech-Sozial Vollekspartei, CSV)+G - The Greens (Déi Gréng, G), known until 1995 as Green List Ecological Alternative-Green Alternative Party (, Gréng Lëscht Ekologesch Initiativ-Gréng Alternativ Partei)'')'),
)

foreach ($ws in $wb.Worksheets) {
    foreach ($pair in $replacements) {
        $null = $ws.Cells.Replace($pair[0], $pair[1], 1, 1, $true, $false, $false, $false)
    }
}